# transactions.xlsx: refresh header labels/casing, drop the "Groceries" row,
# add a new "Laptop" purchase row, renumber the ID column, and store the
# Date column as plain yyyy-mm-dd text instead of Excel date-serial numbers.
#
# Every data cell in the sheet ends up on the sheet's base style (no font /
# number-format override), so every touched cell also gets ClearFormats()'d.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, $text) {
    # Force a plain-text (shared-string) cell with no special number format
    # left behind - matches the "type s / base style" cells in the target
    # sheet. NumberFormat "@" keeps Excel from re-parsing a date-shaped
    # string back into a date serial; the trailing ClearFormats() then
    # drops that temporary Text number-format override again.
    $rng.ClearContents()
    $rng.ClearFormats()
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# --- Header row: description/amount/date -> Description/Amount/Date ----
$ws.Range("B1").ClearFormats()
$ws.Range("B1").Value = "Description"

$ws.Range("C1").ClearFormats()
$ws.Range("C1").Value = "Amount"

$ws.Range("D1").ClearFormats()
$ws.Range("D1").Value = "Date"

# --- Row 2: Rent - amount unchanged, date becomes text ------------------
# A2 (ID=1) and B2 ("Rent") already hold the right values, so they're left
# alone bar the formatting reset.
$ws.Range("B2").ClearFormats()
$ws.Range("C2").ClearFormats()

Set-TextCell $ws.Range("D2") "2025-07-01"

# --- Row 3: Groceries row is dropped; Miscellaneous moves up to ID 3 ----
$ws.Range("A3").ClearFormats()
$ws.Range("A3").Value = 3

$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "Miscellaneous"

$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = 650

Set-TextCell $ws.Range("D3") "2025-07-03"

# --- Row 4: brand-new Laptop purchase, ID 4 ------------------------------
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value = 4

$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "Laptop"

$ws.Range("C4").ClearFormats()
$ws.Range("C4").Value = 1300

Set-TextCell $ws.Range("D4") "2025-07-10"

# --- Column widths (best-fit to the new content) -------------------------
# Target widths are 12.5390625 / 7.515625 / 10.0390625 "characters" (the
# usual 256ths-of-a-character Excel storage). This host quantizes
# ColumnWidth to 1/6-character steps, so these inputs land on the closest
# reachable widths (12.5 / 7.5 / 10.0) instead of matching exactly.
$ws.Columns("B").ColumnWidth = 11.666666666666666
$ws.Columns("C").ColumnWidth = 6.666666666666667
$ws.Columns("D").ColumnWidth = 9.166666666666666
